$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H138").Value = 3391.2932
$ws.Range("I138").Value = 1158.0454
$ws.Range("J138").Value = 4756.0557
$ws.Range("K138").Value = 3474.1362
$ws.Range("L138").Value = 14268.1671
$ws.Range("M138").Value = 1665.8638
$ws.Range("N138").Value = -24548.1671

$ws.Range("H141").Value = 3625.111
$ws.Range("I141").Value = 3232.7144
$ws.Range("J141").Value = 4998.5
$ws.Range("K141").Value = 9698.143199999999
$ws.Range("L141").Value = 14995.5
$ws.Range("M141").Value = -4518.143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4576.757
$ws.Range("I61").Value = 4658.5
$ws.Range("J61").Value = 4514.476
$ws.Range("K61").Value = 4658.5
$ws.Range("L61").Value = 4514.476
$ws.Range("M61").Value = -4446.5
$ws.Range("N61").Value = -4938.476

$ws.Range("H74").Value = 2300.75
$ws.Range("I74").Value = 1328.909
$ws.Range("J74").Value = 3488.5557
$ws.Range("K74").Value = 1328.909
$ws.Range("L74").Value = 3488.5557
$ws.Range("M74").Value = -454.9090000000001

$ws.Range("H77").Value = 2300.75
$ws.Range("I77").Value = 1328.909
$ws.Range("J77").Value = 3488.5557
$ws.Range("K77").Value = 6644.545
$ws.Range("L77").Value = 17442.7785
$ws.Range("M77").Value = -2276.545

$ws.Range("H110").Value = 4498.857
$ws.Range("I110").Value = 2915.5
$ws.Range("J110").Value = 5686.375
$ws.Range("K110").Value = 2915.5
$ws.Range("L110").Value = 5686.375
$ws.Range("M110").Value = -870.5

$ws.Range("H132").Value = 1136548.2
$ws.Range("I132").Value = 1789050.5
$ws.Range("J132").Value = 92544.5
$ws.Range("K132").Value = 5367151.5
$ws.Range("L132").Value = 277633.5
$ws.Range("M132").Value = -5364621.5

$ws.Range("H136").Value = 4576.757
$ws.Range("I136").Value = 4658.5
$ws.Range("J136").Value = 4514.476
$ws.Range("K136").Value = 13975.5
$ws.Range("L136").Value = 13543.428
$ws.Range("M136").Value = -11425.5
$ws.Range("N136").Value = -18643.428

$ws.Range("H137").Value = 130000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 130000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 130000
$ws.Range("N137").Value = -140200

$ws.Range("H139").Value = 74013.8
$ws.Range("I139").Value = 30000
$ws.Range("J139").Value = 85017.25
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 85017.25
$ws.Range("M139").Value = -24860
$ws.Range("N139").Value = -95297.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H95").Value = 12121
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 12121
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 12121
$ws.Range("N95").Value = -17613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2493.3928
$ws.Range("I31").Value = 997.9048
$ws.Range("J31").Value = 3390.6858
$ws.Range("K31").Value = 997.9048
$ws.Range("L31").Value = 3390.6858
$ws.Range("M31").Value = -702.9048
$ws.Range("N31").Value = -3980.6858

$ws.Range("H34").Value = 2493.3928
$ws.Range("I34").Value = 997.9048
$ws.Range("J34").Value = 3390.6858
$ws.Range("K34").Value = 997.9048
$ws.Range("L34").Value = 3390.6858
$ws.Range("M34").Value = -795.9048
$ws.Range("N34").Value = -3794.6858

$ws.Range("H88").Value = 49950
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 49950
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 49950
$ws.Range("N88").Value = -50762

$ws.Range("H91").Value = 49950
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 49950
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 49950
$ws.Range("N91").Value = -52758

$ws.Range("H92").Value = 39999
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 39999
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 39999
$ws.Range("N92").Value = -44991

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 31044.191
$ws.Range("I132").Value = 36721.523
$ws.Range("J132").Value = 7199.4
$ws.Range("K132").Value = 110164.569
$ws.Range("L132").Value = 21598.2
$ws.Range("M132").Value = -107634.569

$ws.Range("H134").Value = 55561470
$ws.Range("I134").Value = 66668268
$ws.Range("J134").Value = 27500
$ws.Range("K134").Value = 200004804
$ws.Range("L134").Value = 82500
$ws.Range("M134").Value = -200002269

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 285.08334
$ws.Range("I60").Value = 306.625
$ws.Range("J60").Value = 242
$ws.Range("K60").Value = 919.875
$ws.Range("L60").Value = 726
$ws.Range("M60").Value = -668.875
$ws.Range("N60").Value = -1228

$ws.Range("H80").Value = 2892.5
$ws.Range("I80").Value = 2890
$ws.Range("J80").Value = 2895
$ws.Range("K80").Value = 8670
$ws.Range("L80").Value = 8685
$ws.Range("M80").Value = -7734
$ws.Range("N80").Value = -10557

$ws.Range("H83").Value = 2892.5
$ws.Range("I83").Value = 2890
$ws.Range("J83").Value = 2895
$ws.Range("K83").Value = 26010
$ws.Range("L83").Value = 26055
$ws.Range("M83").Value = -21330
$ws.Range("N83").Value = -35415

$ws.Range("H87").Value = 22362.375
$ws.Range("I87").Value = 12633
$ws.Range("J87").Value = 28200
$ws.Range("K87").Value = 37899
$ws.Range("L87").Value = 84600
$ws.Range("M87").Value = -36651
$ws.Range("N87").Value = -87096

$ws.Range("H90").Value = 22362.375
$ws.Range("I90").Value = 12633
$ws.Range("J90").Value = 28200
$ws.Range("K90").Value = 113697
$ws.Range("L90").Value = 253800
$ws.Range("M90").Value = -107457
$ws.Range("N90").Value = -266280

$ws.Range("H93").Value = 10000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 30000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -33744

$ws.Range("H137").Value = 9518.5
$ws.Range("I137").Value = 7986.625
$ws.Range("J137").Value = 12582.25
$ws.Range("K137").Value = 23959.875
$ws.Range("L137").Value = 37746.75
$ws.Range("M137").Value = -18859.875
$ws.Range("N137").Value = -47946.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 9820.25
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 9820.25
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 9820.25
$ws.Range("N98").Value = -15810.25

$ws.Range("H113").Value = 6869.037
$ws.Range("I113").Value = 4932.5
$ws.Range("J113").Value = 10742.111
$ws.Range("K113").Value = 4932.5
$ws.Range("L113").Value = 10742.111
$ws.Range("M113").Value = -2762.5
$ws.Range("N113").Value = -15082.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 849
$ws.Range("I22").Value = 849
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 849
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -554
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 849
$ws.Range("I27").Value = 849
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 849
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -742
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 13890201
$ws.Range("I46").Value = 999.1
$ws.Range("J46").Value = 31251702
$ws.Range("K46").Value = 999.1
$ws.Range("L46").Value = 31251702
$ws.Range("M46").Value = -811.1
$ws.Range("N46").Value = -31252078

$ws.Range("H55").Value = 4399.9287
$ws.Range("I55").Value = 1266.3334
$ws.Range("J55").Value = 6750.125
$ws.Range("K55").Value = 1266.3334
$ws.Range("L55").Value = 6750.125
$ws.Range("M55").Value = -1093.3334
$ws.Range("N55").Value = -7096.125

$ws.Range("H93").Value = 7750.125
$ws.Range("I93").Value = 7750.125
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 7750.125
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -6502.125

$ws.Range("H132").Value = 5413
$ws.Range("I132").Value = 3999
$ws.Range("J132").Value = 7298.3335
$ws.Range("K132").Value = 11997
$ws.Range("L132").Value = 21895.0005
$ws.Range("M132").Value = -9467
$ws.Range("N132").Value = -26955.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H109").Value = 35342
$ws.Range("I109").Value = 35342
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 35342
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -33955

$ws.Range("H113").Value = 9264503
$ws.Range("I113").Value = 23811480
$ws.Range("J113").Value = 7335.1816
$ws.Range("K113").Value = 71434440
$ws.Range("L113").Value = 22005.5448
$ws.Range("M113").Value = -71432270
$ws.Range("N113").Value = -26345.5448

$ws.Range("H122").Value = 3392.7334
$ws.Range("I122").Value = 3371.7778
$ws.Range("J122").Value = 3424.1667
$ws.Range("K122").Value = 10115.3334
$ws.Range("L122").Value = 10272.5001
$ws.Range("M122").Value = -7665.3334

$ws.Range("H132").Value = 10294.154
$ws.Range("I132").Value = 8256.333000000001
$ws.Range("J132").Value = 34748
$ws.Range("K132").Value = 24768.999
$ws.Range("L132").Value = 104244
$ws.Range("M132").Value = -22238.999
$ws.Range("N132").Value = -109304
